{"js": "// Load all paragraphs in the document body.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst oldEnding = \"As reuni\u00f5es ser\u00e3o feitas preferencialmente de 10 em 10 dias, de forma online atrav\u00e9s da plataforma do Discord.\";\nconst newEnding = \"As reuni\u00f5es ser\u00e3o feitas preferencialmente de 10 em 10 dias, de forma online atrav\u00e9s da plataforma do Discord;\";\nconst removedText = \"Reuni\u00f5es presenciais ocorrem ocasionalmente e devido a pandemia n\u00e3o s\u00e3o obrigat\u00f3rias.\";\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text.trim();\n\n  if (text === oldEnding) {\n    // Replace the trailing period with a semicolon by rewriting the paragraph text.\n    para.insertText(newEnding, Word.InsertLocation.replace);\n  } else if (text === removedText) {\n    // Remove the entire paragraph (including its paragraph mark).\n    para.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$oldEnding = \"atrav\u00e9s da plataforma do Discord.\"\n$newEnding = \"atrav\u00e9s da plataforma do Discord;\"\n$removedStart = \"Reuni\u00f5es presenciais ocorrem ocasionalmente\"\n\n# Pass 1: drop the whole \"Reuni\u00f5es presenciais ...\" bullet, paragraph mark\n# included. Snapshot + walk back-to-front so a deletion never invalidates\n# the Range of a paragraph we haven't visited yet.\n$paras = @($d.Paragraphs)\nfor ($i = $paras.Count - 1; $i -ge 0; $i--) {\n    if ($paras[$i].Range.Text.StartsWith($removedStart)) {\n        $paras[$i].Range.Delete()\n    }\n}\n\n# Pass 2: turn the trailing period on the Discord bullet into a semicolon.\n# Re-fetch the collection (ranges reseated after pass 1's structural edit),\n# and strip the paragraph-mark character from Range.Text before editing so\n# writing it back doesn't insert a duplicate paragraph break.\nforeach ($para in $d.Paragraphs) {\n    $text = $para.Range.Text\n    if ($text.Contains($oldEnding)) {\n        $mark = $text.Substring($text.Length - 1)\n        $body = $text.Substring(0, $text.Length - 1)\n        $para.Range.Text = $body.Replace($oldEnding, $newEnding)\n    }\n}\n"}
